$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (last changed) date, stored as a date serial.
# Every data row that currently has the old date value (45171) needs to be
# bumped by one day to 45172. Walk rows starting at row 2 (first data row,
# row 1 is the header) until we hit a row whose key column (A) is empty.
$r = 2
while ($true) {
    $keyCell = $ws.Cells.Item($r, 1)
    $key = $keyCell.Value2
    if ([string]::IsNullOrEmpty($key)) {
        break
    }

    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }

    $r = $r + 1
}
